$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: updated instruction text, now wrapped (bold font / formatting retained)
$ws.Range("A1").Value = "Le nom du profil (le nom du rôle ) auquel on veut attribuer les nouveau rôles"
$ws.Range("A1").WrapText = $true

# Move the active selection from A2 to B4
$ws.Range("B4").Select()
